$d = $word.ActiveDocument

$replacements = @(
    @("514÷8=64, 2", "869÷5=173, 4"),
    @("277÷6=46, 1", "562÷9=62, 4"),
    @("299÷2=149, 1", "944÷2=472, 0"),
    @("522÷7=74, 4", "377÷4=94, 1"),
    @("853÷5=170, 3", "303÷7=43, 2"),
    @("786÷2=393, 0", "653÷5=130, 3"),
    @("355÷5=71, 0", "823÷9=91, 4"),
    @("200÷7=28, 4", "577÷5=115, 2"),
    @("889÷5=177, 4", "195÷7=27, 6"),
    @("285÷9=31, 6", "437÷4=109, 1"),
    @("755÷2=377, 1", "195÷8=24, 3"),
    @("647÷9=71, 8", "244÷4=61, 0"),
    @("546÷8=68, 2", "861÷3=287, 0"),
    @("476÷5=95, 1", "455÷5=91, 0"),
    @("243÷7=34, 5", "273÷2=136, 1"),
    @("736÷6=122, 4", "679÷3=226, 1"),
    @("437÷6=72, 5", "315÷8=39, 3"),
    @("661÷3=220, 1", "757÷3=252, 1"),
    @("486÷8=60, 6", "991÷9=110, 1"),
    @("138÷4=34, 2", "724÷3=241, 1"),
    @("388÷8=48, 4", "694÷8=86, 6"),
    @("648÷7=92, 4", "414÷5=82, 4"),
    @("850÷8=106, 2", "662÷8=82, 6"),
    @("329÷9=36, 5", "292÷3=97, 1"),
    @("677÷2=338, 1", "524÷3=174, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
